$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update revised confirmed-death counts for existing rows (per updated source data) ---
$ws.Range("B264").Value = 39
$ws.Range("B270").Value = 43
$ws.Range("B271").Value = 47
$ws.Range("B272").Value = 41
$ws.Range("B273").Value = 43

# --- Append new day row (row 274) ---
$ws.Range("A274").Value = 44172
$ws.Range("A273").Copy()
$ws.Range("A274").PasteSpecial(-4122)

$ws.Range("B274").Value = 45
$ws.Range("D274").Value = 1

$ws.Range("C274").Formula = "=B274+C273"
$ws.Range("E274").Formula = "=D274+E273"
$ws.Range("F274").Formula = "=AVERAGE(B268:B274)"

$excel.CutCopyMode = 0

# --- Restore the active cell/selection on the (now larger) data range ---
$ws.Range("D291").Select() | Out-Null
